$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell W1: "N. Hulkenberg" (same style as the other headers) ---
$ws.Range("W1").Value = "N. Hulkenberg"
$ws.Range("V1").Copy()
$ws.Range("W1").PasteSpecial(-4122)  # xlPasteFormats

# --- Column W gets an (empty) placeholder cell for every existing data row,
#     mirroring the blank "K. Raikkonen" (V) column cells already present. ---
for ($r = 2; $r -le 25; $r++) {
    $cell = $ws.Cells.Item($r, 23)
    $cell.Value = "'"
    $cell.Style = "Normal"
}

# --- Row 26: a brand-new row of driver prices ---
$ws.Range("A26").Value = 44046.91666666666
$ws.Range("A26").NumberFormat = $ws.Range("A25").NumberFormat

$ws.Range("B26").Value = 31.3
$ws.Range("C26").Value = 29.7
$ws.Range("D26").Value = 26
$ws.Range("E26").Value = 23.6
$ws.Range("F26").Value = 20.8
$ws.Range("G26").Value = 20.9
$ws.Range("H26").Value = 15.4
$ws.Range("I26").Value = 13.7
$ws.Range("J26").Value = 12.9
$ws.Range("K26").Value = 12.1
$ws.Range("L26").Value = 10.2

$mcell = $ws.Cells.Item(26, 13)
$mcell.Value = "'"
$mcell.Style = "Normal"
$ncell = $ws.Cells.Item(26, 14)
$ncell.Value = "'"
$ncell.Style = "Normal"

$ws.Range("O26").Value = 9.6
$ws.Range("P26").Value = 10.3
$ws.Range("Q26").Value = 8.800000000000001
$ws.Range("R26").Value = 7.8
$ws.Range("S26").Value = 5.9
$ws.Range("T26").Value = 6.3
$ws.Range("U26").Value = 5.7
$ws.Range("V26").Value = 9.6
$ws.Range("W26").Value = 9.800000000000001
